$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4001.5
$ws.Range("J40").Value = 5002
$ws.Range("L40").Value = 5002
$ws.Range("N40").Value = -5352
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 216.57143
$ws.Range("I5").Value = 216.57143
$ws.Range("K5").Value = 216.57143
$ws.Range("M5").Value = -104.57143
$ws.Range("H32").Value = 22996.475
$ws.Range("I32").Value = 22826.416
$ws.Range("K32").Value = 22826.416
$ws.Range("M32").Value = -22539.416
$ws.Range("H61").Value = 8439.166999999999
$ws.Range("I61").Value = 6533.5415
$ws.Range("K61").Value = 6533.5415
$ws.Range("M61").Value = -6321.5415
$ws.Range("H74").Value = 2835.8635
$ws.Range("I74").Value = 966.2727
$ws.Range("J74").Value = 4705.4546
$ws.Range("K74").Value = 966.2727
$ws.Range("L74").Value = 4705.4546
$ws.Range("M74").Value = -92.27269999999999
$ws.Range("N74").Value = -6453.4546
$ws.Range("H77").Value = 2835.8635
$ws.Range("I77").Value = 966.2727
$ws.Range("J77").Value = 4705.4546
$ws.Range("K77").Value = 4831.363499999999
$ws.Range("L77").Value = 23527.273
$ws.Range("M77").Value = -463.3634999999995
$ws.Range("N77").Value = -32263.273
$ws.Range("H97").Value = 1164654.8
$ws.Range("I97").Value = 1691413.4
$ws.Range("J97").Value = 5785.8
$ws.Range("K97").Value = 1691413.4
$ws.Range("L97").Value = 5785.8
$ws.Range("M97").Value = -1690917.4
$ws.Range("N97").Value = -6777.8
$ws.Range("H122").Value = 3982.6667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3982.6667
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 11948.0001
$ws.Range("N122").Value = -16848.0001
$ws.Range("H132").Value = 5206.6597
$ws.Range("I132").Value = 3651.2942
$ws.Range("K132").Value = 10953.8826
$ws.Range("M132").Value = -8423.882599999999
$ws.Range("H136").Value = 8439.166999999999
$ws.Range("I136").Value = 6533.5415
$ws.Range("K136").Value = 19600.6245
$ws.Range("M136").Value = -17050.6245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 216.57143
$ws.Range("I4").Value = 216.57143
$ws.Range("K4").Value = 216.57143
$ws.Range("M4").Value = -101.57143
$ws.Range("H134").Value = 6685.364
$ws.Range("I134").Value = 4831.273
$ws.Range("K134").Value = 14493.819
$ws.Range("M134").Value = -11958.819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40005280
$ws.Range("I31").Value = 111114010
$ws.Range("K31").Value = 111114010
$ws.Range("M31").Value = -111113715
$ws.Range("H34").Value = 40005280
$ws.Range("I34").Value = 111114010
$ws.Range("K34").Value = 111114010
$ws.Range("M34").Value = -111113808
$ws.Range("H58").Value = 5924.3125
$ws.Range("I58").Value = 5155.8335
$ws.Range("K58").Value = 5155.8335
$ws.Range("M58").Value = -4952.8335
$ws.Range("H62").Value = 13748.75
$ws.Range("I62").Value = 4998.3335
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 4998.3335
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = -4374.3335
$ws.Range("N62").Value = -41248
$ws.Range("H65").Value = 13748.75
$ws.Range("I65").Value = 4998.3335
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 24991.6675
$ws.Range("L65").Value = 200000
$ws.Range("M65").Value = -21871.6675
$ws.Range("N65").Value = -206240
$ws.Range("H132").Value = 193978
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H134").Value = 13306.25
$ws.Range("I134").Value = 14380
$ws.Range("K134").Value = 43140
$ws.Range("M134").Value = -40605
$ws.Range("H136").Value = 5924.3125
$ws.Range("I136").Value = 5155.8335
$ws.Range("K136").Value = 15467.5005
$ws.Range("M136").Value = -12917.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1034.1538
$ws.Range("I140").Value = 873
$ws.Range("K140").Value = 2619
$ws.Range("M140").Value = 2561

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1353.7778
$ws.Range("I122").Value = 1463.4
$ws.Range("J122").Value = 1216.75
$ws.Range("K122").Value = 4390.200000000001
$ws.Range("L122").Value = 3650.25
$ws.Range("M122").Value = -1940.200000000001
$ws.Range("N122").Value = -8550.25
$ws.Range("H132").Value = 6145.8423
$ws.Range("I132").Value = 3042.9092
$ws.Range("K132").Value = 9128.7276
$ws.Range("M132").Value = -6598.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6687.7666
$ws.Range("I46").Value = 1667
$ws.Range("J46").Value = 7245.6294
$ws.Range("K46").Value = 1667
$ws.Range("L46").Value = 7245.6294
$ws.Range("M46").Value = -1479
$ws.Range("N46").Value = -7621.6294
$ws.Range("H100").Value = 11366999
$ws.Range("I100").Value = 35717090
$ws.Range("J100").Value = 3625.0667
$ws.Range("K100").Value = 35717090
$ws.Range("L100").Value = 3625.0667
$ws.Range("M100").Value = -35716549
$ws.Range("N100").Value = -4707.066699999999
$ws.Range("H122").Value = 2829.3333
$ws.Range("I122").Value = 2931.75
$ws.Range("K122").Value = 8795.25
$ws.Range("M122").Value = -6345.25
$ws.Range("H132").Value = 9098.143
$ws.Range("I132").Value = 5801
$ws.Range("J132").Value = 13494.333
$ws.Range("K132").Value = 17403
$ws.Range("L132").Value = 40482.999
$ws.Range("M132").Value = -14873
$ws.Range("N132").Value = -45542.999
$ws.Range("H136").Value = 4776.431
$ws.Range("I136").Value = 4239.108
$ws.Range("J136").Value = 6196.5
$ws.Range("K136").Value = 12717.324
$ws.Range("L136").Value = 18589.5
$ws.Range("M136").Value = -10167.324
$ws.Range("N136").Value = -23689.5
$ws.Range("H140").Value = 112496.75
$ws.Range("J140").Value = 112496.75
$ws.Range("L140").Value = 112496.75
$ws.Range("N140").Value = -122856.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 983.61536
$ws.Range("I113").Value = 1051.65
$ws.Range("J113").Value = 756.8333
$ws.Range("K113").Value = 3154.95
$ws.Range("L113").Value = 2270.4999
$ws.Range("M113").Value = -984.9500000000003
$ws.Range("N113").Value = -6610.4999
$ws.Range("H122").Value = 3411.5
$ws.Range("I122").Value = 3413.1428
$ws.Range("K122").Value = 10239.4284
$ws.Range("M122").Value = -7789.428400000001
$ws.Range("H126").Value = 35842
$ws.Range("I126").Value = 46032.082
$ws.Range("K126").Value = 138096.246
$ws.Range("M126").Value = -135626.246
$ws.Range("H132").Value = 7581.533
$ws.Range("I132").Value = 5463.6665
$ws.Range("J132").Value = 10758.333
$ws.Range("K132").Value = 16390.9995
$ws.Range("L132").Value = 32274.999
$ws.Range("M132").Value = -13860.9995
$ws.Range("N132").Value = -37334.999
$ws.Range("H136").Value = 3166.4482
$ws.Range("I136").Value = 1246.0952
$ws.Range("J136").Value = 8207.375
$ws.Range("K136").Value = 3738.2856
$ws.Range("L136").Value = 24622.125
$ws.Range("M136").Value = -1188.2856
$ws.Range("N136").Value = -29722.125
